$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 58.77403848701053
$ws.Range("C2").Value = 46.90170917246077
$ws.Range("D2").Value = 42.64155957433913
$ws.Range("E2").Value = 41.43963654836019
$ws.Range("F2").Value = 39.75694444444444
$ws.Range("G2").Value = 39.48985007074144
$ws.Range("H2").Value = 39.48985007074144
$ws.Range("I2").Value = 39.40972222222222
